$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append new log rows (23-26, 28-31, 33-36) ---

# Row 23
$ws.Range("A2").Copy()
$ws.Range("A23").PasteSpecial(-4122)
$ws.Range("Z1").Formula = "=""11-25-2023"""
$ws.Range("Z1").Copy()
$ws.Range("A23").PasteSpecial(-4163)
$ws.Range("B23").Value = 'Sreeharsha'
$ws.Range("C2").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value = 'created sendmessage page and styling'

# Row 24
$ws.Range("A2").Copy()
$ws.Range("A24").PasteSpecial(-4122)
$ws.Range("Z1").Formula = "=""11-25-2023"""
$ws.Range("Z1").Copy()
$ws.Range("A24").PasteSpecial(-4163)
$ws.Range("B24").Value = 'Mani Krishna'
$ws.Range("C2").Copy()
$ws.Range("C24").PasteSpecial(-4122)
$ws.Range("C24").Value = 'created book a table page and styling'

# Row 25
$ws.Range("A2").Copy()
$ws.Range("A25").PasteSpecial(-4122)
$ws.Range("Z1").Formula = "=""11-25-2023"""
$ws.Range("Z1").Copy()
$ws.Range("A25").PasteSpecial(-4163)
$ws.Range("B25").Value = 'Jhansi Saranu'
$ws.Range("C2").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$ws.Range("C25").Value = 'created confirmation page and styling'

# Row 26
$ws.Range("A2").Copy()
$ws.Range("A26").PasteSpecial(-4122)
$ws.Range("Z1").Formula = "=""11-25-2023"""
$ws.Range("Z1").Copy()
$ws.Range("A26").PasteSpecial(-4163)
$ws.Range("B26").Value = 'Sunitha Thota'
$ws.Range("C2").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("C26").Value = 'created common css file'

# Row 28
$ws.Range("A2").Copy()
$ws.Range("A28").PasteSpecial(-4122)
$ws.Range("Z1").Formula = "=""12-01-2023"""
$ws.Range("Z1").Copy()
$ws.Range("A28").PasteSpecial(-4163)
$ws.Range("B28").Value = 'Sreeharsha'
$ws.Range("C2").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 'updated styling and responsivenss to the Home page'
$ws.Rows(28).RowHeight = 28.8

# Row 29
$ws.Range("A2").Copy()
$ws.Range("A29").PasteSpecial(-4122)
$ws.Range("Z1").Formula = "=""12-01-2023"""
$ws.Range("Z1").Copy()
$ws.Range("A29").PasteSpecial(-4163)
$ws.Range("B29").Value = 'Mani Krishna'
$ws.Range("C2").Copy()
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C29").Value = 'updated styling and responsivenss to the About us page'
$ws.Rows(29).RowHeight = 28.8

# Row 30
$ws.Range("A2").Copy()
$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("Z1").Formula = "=""12-01-2023"""
$ws.Range("Z1").Copy()
$ws.Range("A30").PasteSpecial(-4163)
$ws.Range("B30").Value = 'Jhansi Saranu'
$ws.Range("C2").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C30").Value = 'updated styling and responsivenss to the Menu page'
$ws.Rows(30).RowHeight = 28.8

# Row 31
$ws.Range("A2").Copy()
$ws.Range("A31").PasteSpecial(-4122)
$ws.Range("Z1").Formula = "=""12-01-2023"""
$ws.Range("Z1").Copy()
$ws.Range("A31").PasteSpecial(-4163)
$ws.Range("B31").Value = 'Sunitha Thota'
$ws.Range("C2").Copy()
$ws.Range("C31").PasteSpecial(-4122)
$ws.Range("C31").Value = 'updated styling and responsivenss to the Contact page'
$ws.Rows(31).RowHeight = 28.8

# Row 33
$ws.Range("A2").Copy()
$ws.Range("A33").PasteSpecial(-4122)
$ws.Range("Z1").Formula = "=""12-03-2023"""
$ws.Range("Z1").Copy()
$ws.Range("A33").PasteSpecial(-4163)
$ws.Range("B33").Value = 'Sreeharsha'
$ws.Range("C2").Copy()
$ws.Range("C33").PasteSpecial(-4122)
$ws.Range("C33").Value = 'created documentation and finalising the home page'
$ws.Rows(33).RowHeight = 28.8

# Row 34
$ws.Range("A2").Copy()
$ws.Range("A34").PasteSpecial(-4122)
$ws.Range("Z1").Formula = "=""12-03-2023"""
$ws.Range("Z1").Copy()
$ws.Range("A34").PasteSpecial(-4163)
$ws.Range("B34").Value = 'Mani Krishna'
$ws.Range("C2").Copy()
$ws.Range("C34").PasteSpecial(-4122)
$ws.Range("C34").Value = 'created documentation and finalising the About us page'
$ws.Rows(34).RowHeight = 28.8

# Row 35
$ws.Range("A2").Copy()
$ws.Range("A35").PasteSpecial(-4122)
$ws.Range("Z1").Formula = "=""12-03-2023"""
$ws.Range("Z1").Copy()
$ws.Range("A35").PasteSpecial(-4163)
$ws.Range("B35").Value = 'Jhansi Saranu'
$ws.Range("C2").Copy()
$ws.Range("C35").PasteSpecial(-4122)
$ws.Range("C35").Value = 'created documentation and finalising the Menu page'
$ws.Rows(35).RowHeight = 28.8

# Row 36
$ws.Range("A2").Copy()
$ws.Range("A36").PasteSpecial(-4122)
$ws.Range("Z1").Formula = "=""12-03-2023"""
$ws.Range("Z1").Copy()
$ws.Range("A36").PasteSpecial(-4163)
$ws.Range("B36").Value = 'Sunitha Thota'
$ws.Range("C2").Copy()
$ws.Range("C36").PasteSpecial(-4122)
$ws.Range("C36").Value = 'created documentation and finalising the Contact us page'
$ws.Rows(36).RowHeight = 28.8

$ws.Range("Z1").ClearContents()

# --- sheet view: drop the saved scroll position, select F10 ---
$ws.Range("F10").Select()
